$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RadarItem")

# Fix the radar icon codes: A20D0x -> A214 0x (rare tool icon fix)
$ws.Range("C3").Value = "A21401"
$ws.Range("C4").Value = "A21402"
$ws.Range("C5").Value = "A21403"
$ws.Range("C6").Value = "A21401"
$ws.Range("C7").Value = "A21402"
$ws.Range("C8").Value = "A21403"

# Restore the selection state on the RadarItem sheet
$ws.Activate()
$ws.Range("C3:C8").Select()
